# Apply updated dSF ("F") column values as described by the commit:
# "repull data, push all data, mean calculation"
# Only column F (dSF) values change for rows 2-22 and 24-27; row 23 is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    4  = 2
    5  = -4
    6  = 4
    7  = 1
    8  = -3
    9  = -4
    10 = 7
    11 = 1
    12 = -4
    13 = -1
    14 = 3
    15 = 4
    16 = -1
    17 = 6
    18 = -3
    19 = -3
    20 = 2
    21 = -4
    22 = -2
    24 = 1
    25 = -3
    26 = -3
    27 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
